# Andre Warsaw ABS Presentation - "Add files via upload" edit
#
# Content changes reproduced here:
#   1. Slide 1 "TextBox 1": trim the byline textbox down to just the
#      "Andre Warsaw" line (drops the "Data 205" / "Prof Lori Perine" /
#      "CRN 21844" paragraphs). The shape has spAutoFit, so its height
#      shrinks automatically once the extra paragraphs are removed.
#   2. Slide 12 "TextBox 4": fix the GitHub link's username run from
#      "wasraw" to "andrewarsaw" without disturbing the surrounding
#      runs/formatting ("https://github.com/" ... "/DATA-205").

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }

        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text

        if ($full.StartsWith("Andre Warsaw") -and $full.Contains("CRN 21844")) {
            # Collapse the textbox back down to a single "Andre Warsaw" paragraph.
            $tr.Text = "Andre Warsaw"
        }
        elseif ($full.Contains("wasraw") -and -not $full.Contains("andrewarsaw")) {
            # Replace just the "wasraw" substring with "andrewarsaw", preserving
            # the run formatting of the untouched text around it.
            $idx = $full.IndexOf("wasraw")
            $sub = $tr.Characters($idx + 1, "wasraw".Length)
            $sub.Text = "andrewarsaw"
        }
    }
}
